# Update master to output generated at c8c62b6
# Replace the two-digit division problems/answers in the table cells.

$d = $word.ActiveDocument

$replacements = @(
    @("31÷3=10, 1", "98÷7=14, 0"),
    @("53÷6=8, 5", "79÷9=8, 7"),
    @("64÷4=16, 0", "91÷8=11, 3"),
    @("51÷3=17, 0", "82÷4=20, 2"),
    @("70÷8=8, 6", "46÷2=23, 0"),
    @("27÷8=3, 3", "20÷9=2, 2"),
    @("36÷2=18, 0", "93÷3=31, 0"),
    @("74÷6=12, 2", "12÷9=1, 3"),
    @("36÷4=9, 0", "19÷6=3, 1"),
    @("67÷2=33, 1", "70÷3=23, 1"),
    @("23÷5=4, 3", "69÷7=9, 6"),
    @("46÷8=5, 6", "63÷3=21, 0"),
    @("24÷5=4, 4", "95÷4=23, 3"),
    @("54÷8=6, 6", "60÷5=12, 0"),
    @("10÷6=1, 4", "72÷7=10, 2"),
    @("22÷4=5, 2", "98÷8=12, 2"),
    @("87÷5=17, 2", "88÷2=44, 0"),
    @("26÷4=6, 2", "54÷9=6, 0"),
    @("62÷8=7, 6", "24÷7=3, 3"),
    @("28÷3=9, 1", "29÷4=7, 1"),
    @("68÷6=11, 2", "57÷9=6, 3"),
    @("38÷6=6, 2", "71÷9=7, 8"),
    @("69÷2=34, 1", "70÷9=7, 7"),
    @("62÷4=15, 2", "90÷6=15, 0"),
    @("11÷7=1, 4", "83÷9=9, 2")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
